# Generate Report for Handback
# Updates the localization-status report: flips the handoff status to
# "handed back" for both languages, records the handback target/back files
# and datetime per-language, and widens a few columns to fit the new text.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$mdFile1 = "3168d6f9-27e9-4f30-9726-229db29279e2.md"
$mdFile2 = "ac51e8fd-9241-4094-a35e-352563aa2e31.md"
$mdUrl1  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d1a2410f4fe3051c38a2f037b94571977486f98d/e2e/3168d6f9-27e9-4f30-9726-229db29279e2.md"
$mdUrl2  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d1a2410f4fe3051c38a2f037b94571977486f98d/e2e/ac51e8fd-9241-4094-a35e-352563aa2e31.md"

# ---------------------------------------------------------------------------
# Overview sheet: per-language status columns (E = zh-cn, F = de-de)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

$wsOverview.Columns.Item(5).ColumnWidth = 29.1
$wsOverview.Columns.Item(6).ColumnWidth = 29.1

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

$wsZh.Range("I2").Value = $mdFile1
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $mdUrl1, "", "", $mdFile1)
$wsZh.Range("J2").Value = "3168d6f9-27e9-4f30-9726-229db29279e2.dff051fc16af6cba3a057da8aaa7802c21545d43.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-15 16:28:33"

$wsZh.Range("I3").Value = $mdFile2
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $mdUrl2, "", "", $mdFile2)
$wsZh.Range("J3").Value = "ac51e8fd-9241-4094-a35e-352563aa2e31.52dcfb707fe6878e52dbe4c85a2e2b3df1f4cde0.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-15 16:28:33"

$wsZh.Columns.Item(3).ColumnWidth = 29.1
$wsZh.Columns.Item(9).ColumnWidth = 39.1
$wsZh.Columns.Item(10).ColumnWidth = 39.1

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

$wsDe.Range("I2").Value = $mdFile1
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $mdUrl1, "", "", $mdFile1)
$wsDe.Range("J2").Value = "3168d6f9-27e9-4f30-9726-229db29279e2.dff051fc16af6cba3a057da8aaa7802c21545d43.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-15 16:28:40"

$wsDe.Range("I3").Value = $mdFile2
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $mdUrl2, "", "", $mdFile2)
$wsDe.Range("J3").Value = "ac51e8fd-9241-4094-a35e-352563aa2e31.52dcfb707fe6878e52dbe4c85a2e2b3df1f4cde0.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-15 16:28:40"

$wsDe.Columns.Item(3).ColumnWidth = 29.1
$wsDe.Columns.Item(9).ColumnWidth = 39.1
$wsDe.Columns.Item(10).ColumnWidth = 39.1
